# Apply cryptos list update (GitHub Actions data refresh) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.757.70"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "2.082.86"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.70"
$ws.Range("E5").Value = "  +0.50%  "
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.08"
$ws.Range("E7").Value = "  +3.36%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0789"
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("E11").Value = "  +2.81%  "
$ws.Range("D12").Value = "2.387.88"
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.78"
$ws.Range("E13").Value = "  +2.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.29"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.775"
$ws.Range("E15").Value = "  +2.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.33"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "2.088.48"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "37.724.94"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.75"
$ws.Range("E20").Value = "  +1.63%  "
$ws.Range("D21").Value = "0.0₃0849"
$ws.Range("E21").Value = "  +3.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.23"
$ws.Range("E22").Value = "  +0.03%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("E24").Value = "  -0.52%  "
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.57"
$ws.Range("E26").Value = "  +7.48%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.53"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.42"
$ws.Range("E29").Value = "  -1.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.56"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("E31").Value = "  +2.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.75"
$ws.Range("E32").Value = "  +2.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0635"
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.68"
$ws.Range("E34").Value = "  +1.94%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.44"
$ws.Range("E36").Value = "  +1.50%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.20"
$ws.Range("E41").Value = "  +1.46%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.26"
$ws.Range("E42").Value = "  +10.57%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0219"
$ws.Range("E43").Value = "  +2.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.91"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").Value = "1.451.83"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.16"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.07"
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.99"
$ws.Range("E50").Value = "  -0.75%  "
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.13"
$ws.Range("E51").Value = "  +0.64%  "
